# Apply crypto price/volume updates per commit "Updated cryptos list on Thu Aug 15 09:54:24 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "58.390.90"
$ws.Range("E2").Value = "  -4.06%  "
$ws.Range("D3").Value = "2.615.76"
$ws.Range("E3").Value = "  -3.40%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "519.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.567"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("E9").Value = "  -6.97%  "
$ws.Range("E10").Value = "  -2.83%  "
$ws.Range("E11").Value = "  -0.65%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").Value = "3.076.98"
$ws.Range("E13").Value = "  -3.44%  "
$ws.Range("D14").Value = "58.369.61"
$ws.Range("E14").Value = "  -3.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.92"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "2.616.73"
$ws.Range("E17").Value = "  -3.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "336.78"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.39"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.39"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.29"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "64.48"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.415"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.98%  "
$ws.Range("E25").Value = "  -1.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.11"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.55%  "
$ws.Range("D28").Value = "0.0₃0789"
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("E29").Value = "  -4.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.76"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.56"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.09"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("E35").Value = "  -3.84%  "
$ws.Range("E36").Value = "  -4.68%  "
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "36.29"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.63"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.57%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.603"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "268.35"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -4.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("E46").Value = "  -4.88%  "
$ws.Range("E47").Value = "  -1.85%  "
$ws.Range("D48").Value = "2.027.46"
$ws.Range("E48").Value = "  -4.80%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0227"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.66"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.21"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.10%  "
